$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "purpose" column (E2:E24) from "S.GISH" to "fullRNASEQ"
$ws.Range("E2:E24").Value = "fullRNASEQ"

# Update the selection to match the diff (single cell E24 selected)
$ws.Range("E24").Select()

# Enable iterative calculation delta setting
$excel.Iteration = $true
$excel.MaxChange = 0.0001
